$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 404.77777
$ws.Range("I17").Value = 100
$ws.Range("J17").Value = 422.70587
$ws.Range("K17").Value = 300
$ws.Range("L17").Value = 1268.11761
$ws.Range("M17").Value = -132
$ws.Range("N17").Value = -1604.11761

$ws.Range("H29").Value = 5000
$ws.Range("J29").Value = 14000
$ws.Range("L29").Value = 42000
$ws.Range("N29").Value = -42562

$ws.Range("H80").Value = 880
$ws.Range("I80").Value = 750
$ws.Range("J80").Value = 966.6667
$ws.Range("K80").Value = 2250
$ws.Range("L80").Value = 2900.0001
$ws.Range("M80").Value = -1252
$ws.Range("N80").Value = -4896.0001

$ws.Range("H82").Value = 17986.75
$ws.Range("I82").Value = 633.3333
$ws.Range("K82").Value = 1899.9999
$ws.Range("M82").Value = -1493.9999

$ws.Range("H83").Value = 880
$ws.Range("I83").Value = 750
$ws.Range("J83").Value = 966.6667
$ws.Range("K83").Value = 6750
$ws.Range("L83").Value = 8700.0003
$ws.Range("M83").Value = -1758
$ws.Range("N83").Value = -18684.0003

$ws.Range("H85").Value = 17986.75
$ws.Range("I85").Value = 633.3333
$ws.Range("K85").Value = 1899.9999
$ws.Range("M85").Value = -495.9999

$ws.Range("H98").Value = 43478800
$ws.Range("I98").Value = 43478800
$ws.Range("K98").Value = 43478800
$ws.Range("M98").Value = -43477302

$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws.Range("H122").Value = 43478800
$ws.Range("I122").Value = 43478800
$ws.Range("K122").Value = 130436400
$ws.Range("M122").Value = -130433950

$ws.Range("H132").Value = 5144.6665
$ws.Range("I132").Value = 6217.8335
$ws.Range("J132").Value = 2998.3333
$ws.Range("K132").Value = 18653.5005
$ws.Range("L132").Value = 8994.999899999999
$ws.Range("M132").Value = -16123.5005
$ws.Range("N132").Value = -14054.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8929722
$ws.Range("I32").Value = 8929722
$ws.Range("K32").Value = 8929722
$ws.Range("M32").Value = -8929435

$ws.Range("H45").Value = 3932.2
$ws.Range("I45").Value = 3887
$ws.Range("K45").Value = 3887
$ws.Range("M45").Value = -3510

$ws.Range("H74").Value = 11915843
$ws.Range("I74").Value = 20835844
$ws.Range("J74").Value = 22507.445
$ws.Range("K74").Value = 20835844
$ws.Range("L74").Value = 22507.445
$ws.Range("M74").Value = -20834970
$ws.Range("N74").Value = -24255.445

$ws.Range("H77").Value = 11915843
$ws.Range("I77").Value = 20835844
$ws.Range("J77").Value = 22507.445
$ws.Range("K77").Value = 104179220
$ws.Range("L77").Value = 112537.225
$ws.Range("M77").Value = -104174852
$ws.Range("N77").Value = -121273.225

$ws.Range("H114").Value = 103197
$ws.Range("J114").Value = 103197
$ws.Range("L114").Value = 103197
$ws.Range("N114").Value = -111875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1377.4546
$ws.Range("I94").Value = 1632.25
$ws.Range("K94").Value = 1632.25
$ws.Range("M94").Value = -1181.25

$ws.Range("H105").Value = 2212.6086
$ws.Range("I105").Value = 1805
$ws.Range("J105").Value = 2657.2727
$ws.Range("K105").Value = 1805
$ws.Range("L105").Value = 2657.2727
$ws.Range("M105").Value = -58
$ws.Range("N105").Value = -6151.2727

$ws.Range("H107").Value = 917.375
$ws.Range("I107").Value = 892.13043
$ws.Range("K107").Value = 892.13043
$ws.Range("M107").Value = 1027.86957

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 452201.97
$ws.Range("I31").Value = 3587.0435
$ws.Range("J31").Value = 808000
$ws.Range("K31").Value = 3587.0435
$ws.Range("L31").Value = 808000
$ws.Range("M31").Value = -3292.0435
$ws.Range("N31").Value = -808590

$ws.Range("H34").Value = 452201.97
$ws.Range("I34").Value = 3587.0435
$ws.Range("J34").Value = 808000
$ws.Range("K34").Value = 3587.0435
$ws.Range("L34").Value = 808000
$ws.Range("M34").Value = -3385.0435
$ws.Range("N34").Value = -808404

$ws.Range("H105").Value = 3582.625
$ws.Range("I105").Value = 3117.3333
$ws.Range("J105").Value = 3861.8
$ws.Range("K105").Value = 3117.3333
$ws.Range("L105").Value = 3861.8
$ws.Range("M105").Value = -1370.3333
$ws.Range("N105").Value = -7355.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 559094.3
$ws.Range("J12").Value = 950433.4
$ws.Range("L12").Value = 2851300.2
$ws.Range("N12").Value = -2851646.2

$ws.Range("H98").Value = 1000
$ws.Range("I98").Value = 1000
$ws.Range("J98").Value = 1000
$ws.Range("K98").Value = 3000
$ws.Range("L98").Value = 3000
$ws.Range("M98").Value = -1502
$ws.Range("N98").Value = -5996

$ws.Range("H136").Value = 8466.888999999999
$ws.Range("I136").Value = 5771.7144
$ws.Range("K136").Value = 17315.1432
$ws.Range("M136").Value = -12215.1432

$ws.Range("H137").Value = 6314.25
$ws.Range("I137").Value = 5383.1665
$ws.Range("J137").Value = 7245.3335
$ws.Range("K137").Value = 16149.4995
$ws.Range("L137").Value = 21736.0005
$ws.Range("M137").Value = -11049.4995
$ws.Range("N137").Value = -31936.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1861.2778
$ws.Range("I107").Value = 1134.5834
$ws.Range("J107").Value = 3314.6667
$ws.Range("K107").Value = 1134.5834
$ws.Range("L107").Value = 3314.6667
$ws.Range("M107").Value = 785.4166
$ws.Range("N107").Value = -7154.6667

$ws.Range("H113").Value = 4043.8667
$ws.Range("I113").Value = 3725.75
$ws.Range("J113").Value = 5316.3335
$ws.Range("K113").Value = 3725.75
$ws.Range("L113").Value = 5316.3335
$ws.Range("M113").Value = -1555.75
$ws.Range("N113").Value = -9656.333500000001

$ws.Range("H122").Value = 1599.1818
$ws.Range("I122").Value = 1599.1818
$ws.Range("K122").Value = 4797.5454
$ws.Range("M122").Value = -2347.5454

$ws.Range("H132").Value = 32261646
$ws.Range("I132").Value = 35717916
$ws.Range("K132").Value = 107153748
$ws.Range("M132").Value = -107151218

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3839.25
$ws.Range("I16").Value = 3854
$ws.Range("J16").Value = 3818.6
$ws.Range("K16").Value = 3854
$ws.Range("L16").Value = 3818.6
$ws.Range("M16").Value = -3684
$ws.Range("N16").Value = -4158.6

$ws.Range("H55").Value = 62500412
$ws.Range("I55").Value = 90909530
$ws.Range("K55").Value = 90909530
$ws.Range("M55").Value = -90909357

$ws.Range("H68").Value = 2528.6
$ws.Range("I68").Value = 2532
$ws.Range("J68").Value = 2498
$ws.Range("K68").Value = 2532
$ws.Range("L68").Value = 2498
$ws.Range("M68").Value = -1783
$ws.Range("N68").Value = -3996

$ws.Range("H71").Value = 2528.6
$ws.Range("I71").Value = 2532
$ws.Range("J71").Value = 2498
$ws.Range("K71").Value = 12660
$ws.Range("L71").Value = 12490
$ws.Range("M71").Value = -8916
$ws.Range("N71").Value = -19978
